# Start evaluation and cleaning of documents:
# clear the already-coded values in row 14 (this entry had not actually
# been reviewed yet) while leaving the relevance flag (B14) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C14:G14").ClearContents()

# Reflect where the reviewer was working in the sheet view/selection.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G14").Select()

$wb.Save()
